# =====================================================================
# Edit script: implement "assays version 1" change to scrnaseq-metadata
# - Insert a new "version list" sheet (value "1") right after "Export as TSV"
# - Insert two new leading columns in "Export as TSV": "version" and
#   "description", each with its own header comment and (for "version")
#   a list-based data validation against the new "version list" sheet.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$headers = @(
    "version"
    "description"
    "donor_id"
    "tissue_id"
    "execution_datetime"
    "protocols_io_doi"
    "operator"
    "operator_email"
    "pi"
    "pi_email"
    "assay_category"
    "assay_type"
    "analyte_class"
    "is_targeted"
    "acquisition_instrument_vendor"
    "acquisition_instrument_model"
    "sc_isolation_protocols_io_doi"
    "sc_isolation_entity"
    "sc_isolation_tissue_dissociation"
    "sc_isolation_enrichment"
    "sc_isolation_quality_metric"
    "sc_isolation_cell_number"
    "rnaseq_assay_input"
    "rnaseq_assay_method"
    "library_construction_protocols_io_doi"
    "library_layout"
    "library_adapter_sequence"
    "library_id"
    "is_technical_replicate"
    "cell_barcode_read"
    "cell_barcode_offset"
    "cell_barcode_size"
    "library_pcr_cycles"
    "library_pcr_cycles_for_sample_index"
    "library_final_yield_value"
    "library_final_yield_unit"
    "library_average_fragment_size"
    "sequencing_reagent_kit"
    "sequencing_read_format"
    "sequencing_read_percent_q30"
    "sequencing_phix_percent"
    "contributors_path"
    "data_path"
)

$commentTexts = @(
    "Version of the schema to use when validating this metadata."
    "Free-text description of this assay."
    "HuBMAP Display ID of the donor of the assayed tissue."
    "HuBMAP Display ID of the assayed tissue."
    "Start date and time of assay, typically a date-time stamped folder generated by the acquisition instrument. YYYY-MM-DD hh:mm, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s, hh is the hour with leading zeros, mm are the minutes with leading zeros."
    "DOI for protocols.io referring to the protocol for this assay."
    "Name of the person responsible for executing the assay."
    "Email address for the operator."
    "Name of the principal investigator responsible for the data."
    "Email address for the principal investigator."
    "Each assay is placed into one of the following 3 general categories: generation of images of microscopic entities, identification & quantitation of molecules by mass spectrometry, and determination of nucleotide sequence."
    "The specific type of assay being executed."
    "Analytes are the target molecules being measured with the assay."
    "Specifies whether or not a specific molecule(s) is/are targeted for detection/measurement by the assay. The CODEX analyte is protein."
    "An acquisition instrument is the device that contains the signal detection hardware and signal processing software. Assays generate signals such as light of various intensities or color or signals representing the molecular mass."
    "Manufacturers of an acquisition instrument may offer various versions (models) of that instrument with different features or sensitivities. Differences in features or sensitivities may be relevant to processing or interpretation of the data."
    "Link to a protocols document answering the question: How were single cells separated into a single-cell suspension?"
    "The type of single cell entity derived from isolation protocol"
    "The method by which tissues are dissociated into single cells in suspension."
    "The method by which specific cell populations are sorted or enriched."
    "A quality metric by visual inspection prior to cell lysis or defined by known parameters such as wells with several cells or no cells. This can be captured at a high level."
    "Total number of cell/nuclei yielded post dissociation and enrichment"
    "Number of cell/nuclei input to the assay"
    "The kit used for the RNA sequencing assay"
    "A link to the protocol document containing the library construction method (including version) that was used, e.g. `"Smart-Seq2`", `"Drop-Seq`", `"10X v3`"."
    "Whether the library was generated for single-end or paired end sequencing"
    "Adapter sequence to be used for adapter trimming"
    "An id for the library. The id may be text and/or numbers"
    "Is the sequencing reaction run in repliucate, TRUE or FALSE"
    "Which read file contains the cell barcode"
    "Position(s) in the read at which the cell barcode starts."
    "Length of the cell barcode in base pairs"
    "Number of PCR cycles to amplify cDNA"
    "Number of PCR cycles performed for library indexing"
    "Total number of ng of library after final pcr amplification step. This is the concentration (ng/ul) * volume (ul)"
    "Units of final library yield"
    "Average size of sequencing library fragments estimated via gel electrophoresis or bioanalyzer/tapestation."
    "Reagent kit used for sequencing"
    "Slash-delimited list of the number of sequencing cycles for, for example, Read1, i7 index, i5 index, and Read2."
    "Percent of bases with Quality scores above Q30"
    "Percent PhiX loaded to the run"
    "Relative path to file with ORCID IDs for contributors for this dataset."
    "Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions."
)

# ---------------------------------------------------------------------
# 1. Remove all existing header comments on sheet 1. The comments in
#    this engine are anchored to a fixed cell position and do NOT shift
#    automatically when columns are inserted, so we clear them first and
#    re-create them afterwards at their correct final positions.
# ---------------------------------------------------------------------
$existingCount = $ws.Comments.Count
for ($i = $existingCount; $i -ge 1; $i--) {
    $ws.Comments.Item($i).Delete()
}

# ---------------------------------------------------------------------
# 2. Insert two new blank columns at the front of the sheet (A and B).
#    Everything that used to be in column A now lives in column C, etc.
#    Data validations defined on the old columns shift automatically.
# ---------------------------------------------------------------------
$ws.Range("A:B").Insert()

# ---------------------------------------------------------------------
# 3. Insert the new "version list" worksheet right after "Export as TSV"
#    and give it its single value, "1", stored as text (matching the
#    original workbook's convention of storing list values as strings).
# ---------------------------------------------------------------------
$versionList = $wb.Worksheets.Add($null, $ws)
$versionList.Name = "version list"
$versionList.Range("A1").Value = "'1"
$versionList.Range("A1").Style = "Normal"

# ---------------------------------------------------------------------
# 4. Rewrite the header row (row 1) of "Export as TSV" with the full,
#    final set of column names (version, description, then all of the
#    pre-existing names in their original relative order).
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}

# The two brand-new header cells (A1, B1) come in with default
# formatting; match them to the bold/centered/wrap-text style already
# used by the rest of the header row.
$newHeaderCells = $ws.Range("A1:B1")
$newHeaderCells.Font.Bold = $true
$newHeaderCells.HorizontalAlignment = -4108
$newHeaderCells.WrapText = $true

# ---------------------------------------------------------------------
# 5. Re-create the header comments in their final positions.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $commentTexts.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.AddComment($commentTexts[$i])
}

# ---------------------------------------------------------------------
# 6. Add the data validation for the new "version" column (A), matching
#    the style of the other list-based validations already present.
# ---------------------------------------------------------------------
$verRange = $ws.Range("A2:A1048576")
$verRange.Validation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1")
$verRange.Validation.ErrorTitle = "Value must come from list"
$verRange.Validation.ErrorMessage = "Value must be one of: 1."

Write-Host "Done."
